# "adding csv and updated modal"
# Populates Sheet2 with the alcohol / solo-cup calculations, switches the
# active sheet from Sheet1 to Sheet2, and updates the selections on both
# sheets to match the post-edit workbook state.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# Sheet2 header row (row 1). Written in the exact order the strings were
# first introduced so new shared-string indices line up with the source
# workbook (C1..G1 first, the row-7 labels next, then H1, then A14).
# ---------------------------------------------------------------------
$ws2.Range("A1").Value = "alcohol type"
$ws2.Range("B1").Value = "alcohol %"
$ws2.Range("C1").Value = "standard drink size (oz)"
$ws2.Range("D1").Value = "standard drink size (L)"
$ws2.Range("E1").Value = "full red cup (oz)"
$ws2.Range("F1").Value = "full red cup (L)"
$ws2.Range("G1").Value = "pure alcohol per standard drink (L)"

$ws2.Range("A7").Value = "Pure Alcohol (L)"
$ws2.Range("G7").Value = "number of red cup standard drinks"
$ws2.Range("H7").Value = "number of full cup drinks"

$ws2.Range("H1").Value = "pure alcohol per full cup (L)"

$ws2.Range("A14").Value = "pure Alcohol (oz)"

# ---------------------------------------------------------------------
# Row 2 - beer
# ---------------------------------------------------------------------
$ws2.Range("A2").Value = "beer"
$ws2.Range("B2").Value = 0.05
$ws2.Range("B2").NumberFormat = "0%"
$ws2.Range("C2").Value = 12
$ws2.Range("D2").Formula = "=C2/33.814"
$ws2.Range("E2").Value = 16
$ws2.Range("F2").Formula = "=E2/33.814"
$ws2.Range("G2").Formula = "=D2 * B2"
$ws2.Range("H2").Formula = "=F2 *B2"
$ws2.Range("K2").Formula = "= 0.01774413 * 33.814"

# ---------------------------------------------------------------------
# Row 3 - wine, Row 4 - spirits
# ---------------------------------------------------------------------
$ws2.Range("A3").Value = "wine"
$ws2.Range("B3").Value = 0.12
$ws2.Range("B3").NumberFormat = "0%"
$ws2.Range("C3").Value = 5
$ws2.Range("E3").Value = 16

$ws2.Range("A4").Value = "spirits"
$ws2.Range("B4").Value = 0.4
$ws2.Range("B4").NumberFormat = "0%"
$ws2.Range("C4").Value = 1.5
$ws2.Range("E4").Value = 16

# Shared formulas spanning rows 3:4
$ws2.Range("D3:D4").Formula = "=C3/33.814"
$ws2.Range("F3:F4").Formula = "=E3/33.814"
$ws2.Range("G3:G4").Formula = "=D3 * B3"
$ws2.Range("H3:H4").Formula = "=F3 *B3"

# ---------------------------------------------------------------------
# Rows 8-10 - standard-drinks-per-bottle summary
# ---------------------------------------------------------------------
$ws2.Range("A8").Value = 3.6

$ws2.Range("F8").Value = "beer"
$ws2.Range("G8").Formula = "=A8/G2"
$ws2.Range("H8").Formula = "=A8/H2"

$ws2.Range("F9").Value = "wine"
$ws2.Range("G9").Formula = "=A8/G3"
$ws2.Range("H9").Formula = "=A8/H3"

$ws2.Range("F10").Value = "spirits"
$ws2.Range("G10").Formula = "=A8/G4"
$ws2.Range("H10").Formula = "=A8/H4"

# ---------------------------------------------------------------------
# Row 15 - pure alcohol (oz) cross-check
# ---------------------------------------------------------------------
$ws2.Range("A15").Formula = "=A8 * 33.814"
$ws2.Range("G15").Formula = "= A15/K2"

# ---------------------------------------------------------------------
# Column widths (best effort - headless engine rounds to the nearest
# 1/6 of a character, so these land as close as possible to the widths
# Excel's own AutoFit produced).
# ---------------------------------------------------------------------
$ws2.Columns.Item(1).ColumnWidth = 12.944010416666666
$ws2.Columns.Item(3).ColumnWidth = 19.053385416666668
$ws2.Columns.Item(4).ColumnWidth = 17.830729166666668
$ws2.Columns.Item(5).ColumnWidth = 15.276041666666666
$ws2.Columns.Item(6).ColumnWidth = 16.053385416666668
$ws2.Columns.Item(7).ColumnWidth = 29.276041666666668
$ws2.Columns.Item(8).ColumnWidth = 24.385416666666668

# ---------------------------------------------------------------------
# Selections + active sheet: Sheet1 keeps a block selection (A1:B4) and
# loses focus; Sheet2 becomes the active/visible tab with G16 selected.
# ---------------------------------------------------------------------
[void]$ws1.Range("A1:B4").Select()
[void]$ws2.Range("G16").Select()
[void]$ws2.Activate()
